$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-07-10 Wednesday"; new = "2024-07-11 Thursday"},
    @{old = "286×8="; new = "340×2="},
    @{old = "273×6="; new = "264×3="},
    @{old = "936×7="; new = "702×3="},
    @{old = "356×8="; new = "767×9="},
    @{old = "203×2="; new = "559×8="},
    @{old = "778×6="; new = "124×8="},
    @{old = "861×2="; new = "821×3="},
    @{old = "916×4="; new = "546×5="},
    @{old = "929×6="; new = "688×7="},
    @{old = "543×9="; new = "233×9="},
    @{old = "564×2="; new = "884×5="},
    @{old = "821×4="; new = "784×2="},
    @{old = "990×9="; new = "512×3="},
    @{old = "748×7="; new = "702×2="},
    @{old = "996×5="; new = "613×7="},
    @{old = "589×6="; new = "463×2="},
    @{old = "241×4="; new = "488×5="},
    @{old = "257×3="; new = "675×6="},
    @{old = "125×3="; new = "141×6="},
    @{old = "896×4="; new = "886×3="},
    @{old = "399×7="; new = "413×3="},
    @{old = "137×3="; new = "544×6="},
    @{old = "895×3="; new = "925×4="},
    @{old = "990×6="; new = "317×3="},
    @{old = "310×3="; new = "295×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
